$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 28: "Num tanks" label (bold) in I28, description in G28 ---
$ws.Range("I28").Value = "Num tanks"
$ws.Range("I28").Font.Bold = $true
$ws.Range("G28").Value = "Alage head tank fed to 8 tanks "

# --- Row 29: number of tanks ---
$ws.Range("I29").Value = 8

# --- Row 30: peristaltic note + flow rate label (bold) ---
$ws.Range("G30").Value = "Peristaltic was 2 mL per minute"
$ws.Range("I30").Value = "Flow rate (mL per minute)"
$ws.Range("I30").Font.Bold = $true

# --- Row 31: flow rate value (mL per minute) ---
$ws.Range("I31").Value = 2

# --- Row 32: flow rate per day note ---
$ws.Range("G32").Value = "Flow rate per day "

# --- Row 33: flow rate per day (mL), label bold, extend bold through L33 ---
$ws.Range("H33").Formula = "=I31*60*24"
$ws.Range("I33").Value = "mL algae per high-food tank per day!"
$ws.Range("H33:L33").Font.Bold = $true

# --- Row 35: High chl cells per day per tank note ---
$ws.Range("G35").Value = "High chl cells per day per tank ==[ (mean high chl cells as cells mL * mL per tank per day) / 8]"

# --- Row 36: High chl cells per day per tank formula ---
$ws.Range("G36").Formula = "=(H22*H33)/8"
$ws.Range("G36").Style = "Normal"

# --- Row 37: Low chl cells per day per tank note ---
$ws.Range("G37").Value = "Low chl cells per day per tank ==[ (mean high chl cells as cells mL * mL per tank per day) / 8]"

# --- Row 38: Low chl cells per day per tank formula ---
$ws.Range("G38").Formula = "=(H25*H33)/8"
$ws.Range("G38").Style = "Normal"

# --- View changes: scroll to D22 and select J27 ---
$win = $ws.Application.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 4
[void]$ws.Range("J27").Select()
